# Update the worksheet date header.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-04-15 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-04-22 Tuesday", 2)

# Update the division-problem table. The table has 20 rows; only rows
# 1, 5, 9, 13 and 17 (1-indexed) carry problem text, the rest are blank
# answer rows. Each of those 5 rows has 5 cells (columns 1-5).
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "29÷2=14, 1"
$t.Cell(1, 2).Range.Text = "35÷2=17, 1"
$t.Cell(1, 3).Range.Text = "36÷9=4, 0"
$t.Cell(1, 4).Range.Text = "90÷9=10, 0"
$t.Cell(1, 5).Range.Text = "53÷5=10, 3"

# Row 5
$t.Cell(5, 1).Range.Text = "70÷6=11, 4"
$t.Cell(5, 2).Range.Text = "34÷3=11, 1"
$t.Cell(5, 3).Range.Text = "32÷4=8, 0"
$t.Cell(5, 4).Range.Text = "60÷6=10, 0"
$t.Cell(5, 5).Range.Text = "43÷7=6, 1"

# Row 9
$t.Cell(9, 1).Range.Text = "20÷5=4, 0"
$t.Cell(9, 2).Range.Text = "11÷6=1, 5"
$t.Cell(9, 3).Range.Text = "57÷4=14, 1"
$t.Cell(9, 4).Range.Text = "85÷4=21, 1"
$t.Cell(9, 5).Range.Text = "75÷4=18, 3"

# Row 13
$t.Cell(13, 1).Range.Text = "94÷8=11, 6"
$t.Cell(13, 2).Range.Text = "88÷6=14, 4"
$t.Cell(13, 3).Range.Text = "34÷5=6, 4"
$t.Cell(13, 4).Range.Text = "71÷9=7, 8"
$t.Cell(13, 5).Range.Text = "49÷8=6, 1"

# Row 17
$t.Cell(17, 1).Range.Text = "83÷3=27, 2"
$t.Cell(17, 2).Range.Text = "86÷4=21, 2"
$t.Cell(17, 3).Range.Text = "37÷3=12, 1"
$t.Cell(17, 4).Range.Text = "87÷8=10, 7"
$t.Cell(17, 5).Range.Text = "29÷7=4, 1"

Write-Host "Done updating document."
